$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Voltmeter")

# Update prices in column B
$ws1.Range("B2").Value = 68
$ws1.Range("B3").Value = 77
$ws1.Range("B4").Value = 120
$ws1.Range("B5").Value = 134
$ws1.Range("B6").Value = 199
$ws1.Range("B7").Value = 61
$ws1.Range("B8").Value = 68
$ws1.Range("B9").Value = 84

# Widen column A on the Voltmeter sheet (target raw width = 44)
$ws1.Columns("A").ColumnWidth = 43.166666666666664

# Move the selection on the Voltmeter sheet
$ws1.Range("C11").Select() | Out-Null

# Add a new "Metadata" sheet after the existing one
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Metadata"

$newSheet.Range("A1").Value = "Энергохит"
$newSheet.Range("B1").Value = "'07.24.2013"
$newSheet.Range("C1").Value = "' 01.08.2012"
$newSheet.Range("D1").Value = "Updated prices"
$newSheet.Range("A1:D1").ClearFormats()

$newSheet.Range("E6").Select() | Out-Null
